# Update the division-fact answers in the single table of the document.
# Each data row in the table is followed by 3 blank spacer rows, so the
# populated rows are Word table rows 1, 5, 9, 13, 17 (1-indexed).
#
# Cell values are addressed by (row, column) rather than via Find/Replace
# so that the two identical "84÷7=12, 0" cells in row 17 can be updated to
# their two different target values without ambiguity.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "64÷6=10, 4" },
    @{ Row = 1;  Col = 2; Text = "27÷5=5, 2" },
    @{ Row = 1;  Col = 3; Text = "83÷6=13, 5" },
    @{ Row = 1;  Col = 4; Text = "43÷8=5, 3" },
    @{ Row = 1;  Col = 5; Text = "22÷9=2, 4" },

    @{ Row = 5;  Col = 1; Text = "97÷8=12, 1" },
    @{ Row = 5;  Col = 2; Text = "38÷6=6, 2" },
    @{ Row = 5;  Col = 3; Text = "46÷8=5, 6" },
    @{ Row = 5;  Col = 4; Text = "82÷8=10, 2" },
    @{ Row = 5;  Col = 5; Text = "52÷9=5, 7" },

    @{ Row = 9;  Col = 1; Text = "21÷3=7, 0" },
    @{ Row = 9;  Col = 2; Text = "34÷9=3, 7" },
    @{ Row = 9;  Col = 3; Text = "74÷4=18, 2" },
    @{ Row = 9;  Col = 4; Text = "44÷6=7, 2" },
    @{ Row = 9;  Col = 5; Text = "51÷8=6, 3" },

    @{ Row = 13; Col = 1; Text = "76÷9=8, 4" },
    @{ Row = 13; Col = 2; Text = "67÷7=9, 4" },
    @{ Row = 13; Col = 3; Text = "96÷7=13, 5" },
    @{ Row = 13; Col = 4; Text = "47÷3=15, 2" },
    @{ Row = 13; Col = 5; Text = "40÷8=5, 0" },

    @{ Row = 17; Col = 1; Text = "76÷3=25, 1" },
    @{ Row = 17; Col = 2; Text = "22÷4=5, 2" },
    @{ Row = 17; Col = 3; Text = "18÷6=3, 0" },
    @{ Row = 17; Col = 4; Text = "57÷8=7, 1" },
    @{ Row = 17; Col = 5; Text = "15÷8=1, 7" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
